# Updates the stats table in the Renaissance / JDK21 / ZGC chi-square
# (heap-1G) benchmark doc (single-column table, 46 rows):
#  - rows 1-3 (summary figures) become "0M"
#  - row 4 count changes 653 -> 869
#  - rows 7-8 percentile figures get refreshed
#  - row 12 total time figure gets refreshed
#  - the three per-run breakdown rows (44-46), each originally a run of
#    tab-separated values, collapse down to the single summary figure
#    that used to live in rows 1-3 (86.63 / 14.65 / 109)

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellText($rowIndex, $text) {
    $cell = $tbl.Cell($rowIndex, 1)
    $cell.Range.Text = $text
}

Set-CellText 1 "0M"
Set-CellText 2 "0M"
Set-CellText 3 "0M"
Set-CellText 4 "869"

Set-CellText 7 "0.05135"
Set-CellText 8 "0.00645"

Set-CellText 12 "14.64527"

Set-CellText 44 "86.63"
Set-CellText 45 "14.65"
Set-CellText 46 "109"
